$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 25.1900000000005
$ws.Range("H2").Value = 0.0000000000000002589441456851677
$ws.Range("K2").Value = 44.21298096435004
$ws.Range("L2").Value = "[33.24697231653519, 55.178989612164884]"
$ws.Range("M2").Value = 0.0000000000000333066907387547
$ws.Range("N2").Value = 0.0000000000000333066907387547
$ws.Range("O2").Value = 1.867974010242579
$ws.Range("P2").Value = "[1.603816069400195, 2.132131951084964]"
$ws.Range("S2").Value = 63.36649422119477
$ws.Range("T2").Value = "[57.17275084257348, 69.56023759981606]"
$ws.Range("W2").Value = 17.70108108108143
$ws.Range("X2").Value = 16.64204204204238
$ws.Range("Y2").Value = 18.76012012012049

# Row 3
$ws.Range("E3").Value = 23.03000000000016
$ws.Range("H3").Value = 0.0000000000000002589441456851677
$ws.Range("K3").Value = 44.92934297195973
$ws.Range("L3").Value = "[34.58542200378896, 55.273263940130505]"
$ws.Range("M3").Value = 0.0000000000000004440892098500626
$ws.Range("N3").Value = 0.0000000000000008881784197001252
$ws.Range("O3").Value = -1.886842434588464
$ws.Range("P3").Value = "[-2.125842476303003, -1.647842392873926]"
$ws.Range("S3").Value = 62.96612689680516
$ws.Range("T3").Value = "[57.245981104150204, 68.68627268946011]"
$ws.Range("W3").Value = 6.915915915915964
$ws.Range("X3").Value = 6.039899899899943
$ws.Range("Y3").Value = 7.791931931931984
